$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jeunes")
$ws.Range("D2").Value = '<jt:if test="${adherent.agecalcule}">${adherent.age}</jt:if>'
$ws.Range("E2").Value = '<jt:if test="${adherent.agecalcule}">${adherent.agecamp}</jt:if>'
